$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1949458483754513
$ws.Range("C2").Value = 0.5595667870036101
$ws.Range("J2").Value = 0.02166064981949458
$ws.Range("P2").Value = 0.1299638989169675
$ws.Range("S2").Value = 0.09386281588447654
$ws.Range("B3").Value = 0.02339181286549707
$ws.Range("C3").Value = 0.07017543859649122
$ws.Range("J3").Value = 0.01754385964912281
$ws.Range("P3").Value = 0.7543859649122807
$ws.Range("S3").Value = 0.1345029239766082
$ws.Range("J4").Value = 0.02325581395348837
$ws.Range("P4").Value = 0.6279069767441861
$ws.Range("S4").Value = 0.3488372093023256
$ws.Range("J5").Value = 0.3333333333333333
$ws.Range("P5").Value = 0.3333333333333333
$ws.Range("S5").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.05429864253393665
$ws.Range("D6").Value = 0.009049773755656109
$ws.Range("F6").Value = 0.07692307692307693
$ws.Range("J6").Value = 0.248868778280543
$ws.Range("O6").Value = 0.03167420814479638
$ws.Range("Q6").Value = 0.1719457013574661
$ws.Range("R6").Value = 0.07239819004524888
$ws.Range("S6").Value = 0.334841628959276
$ws.Range("B7").Value = 0.08294930875576037
$ws.Range("D7").Value = 0.02304147465437788
$ws.Range("E7").Value = 0.004608294930875576
$ws.Range("F7").Value = 0.04608294930875576
$ws.Range("J7").Value = 0.152073732718894
$ws.Range("O7").Value = 0.02304147465437788
$ws.Range("Q7").Value = 0.1981566820276498
$ws.Range("R7").Value = 0.08755760368663594
$ws.Range("S7").Value = 0.3824884792626728
$ws.Range("B8").Value = 0.0796812749003984
$ws.Range("D8").Value = 0.01394422310756972
$ws.Range("E8").Value = 0.00199203187250996
$ws.Range("F8").Value = 0.05577689243027888
$ws.Range("J8").Value = 0.149402390438247
$ws.Range("O8").Value = 0.02390438247011952
$ws.Range("Q8").Value = 0.1533864541832669
$ws.Range("R8").Value = 0.09163346613545817
$ws.Range("S8").Value = 0.4302788844621514
$ws.Range("B9").Value = 0.1043956043956044
$ws.Range("D9").Value = 0.01648351648351648
$ws.Range("F9").Value = 0.07692307692307693
$ws.Range("J9").Value = 0.1483516483516484
$ws.Range("O9").Value = 0.02197802197802198
$ws.Range("Q9").Value = 0.1428571428571428
$ws.Range("R9").Value = 0.1373626373626374
$ws.Range("S9").Value = 0.3516483516483517
$ws.Range("B10").Value = 0.09778456837280367
$ws.Range("D10").Value = 0.01986249045072574
$ws.Range("E10").Value = 0.0007639419404125286
$ws.Range("F10").Value = 0.06417112299465241
$ws.Range("J10").Value = 0.1291061879297173
$ws.Range("O10").Value = 0.01527883880825057
$ws.Range("Q10").Value = 0.2077922077922078
$ws.Range("R10").Value = 0.08785332314744079
$ws.Range("S10").Value = 0.3773873185637892
$ws.Range("G11").Value = 0.1412742382271468
$ws.Range("J11").Value = 0.09418282548476455
$ws.Range("K11").Value = 0.1994459833795014
$ws.Range("L11").Value = 0.5512465373961218
$ws.Range("S11").Value = 0.01385041551246537
$ws.Range("G12").Value = 0.7438423645320197
$ws.Range("J12").Value = 0.167487684729064
$ws.Range("K12").Value = 0.009852216748768473
$ws.Range("L12").Value = 0.03448275862068965
$ws.Range("S12").Value = 0.04433497536945813
$ws.Range("G13").Value = 0.7575757575757576
$ws.Range("J13").Value = 0.2121212121212121
$ws.Range("S13").Value = 0.0303030303030303
$ws.Range("F15").Value = 0.03478260869565217
$ws.Range("H15").Value = 0.191304347826087
$ws.Range("I15").Value = 0.06521739130434782
$ws.Range("J15").Value = 0.3260869565217391
$ws.Range("K15").Value = 0.09130434782608696
$ws.Range("M15").Value = 0.01304347826086956
$ws.Range("N15").Value = 0.004347826086956522
$ws.Range("O15").Value = 0.06521739130434782
$ws.Range("S15").Value = 0.208695652173913
$ws.Range("F16").Value = 0.02702702702702703
$ws.Range("H16").Value = 0.1945945945945946
$ws.Range("I16").Value = 0.03783783783783784
$ws.Range("J16").Value = 0.4162162162162162
$ws.Range("K16").Value = 0.1189189189189189
$ws.Range("M16").Value = 0.03783783783783784
$ws.Range("N16").Value = 0.005405405405405406
$ws.Range("O16").Value = 0.04324324324324325
$ws.Range("S16").Value = 0.1189189189189189
$ws.Range("F17").Value = 0.01535087719298246
$ws.Range("H17").Value = 0.1973684210526316
$ws.Range("I17").Value = 0.09649122807017543
$ws.Range("J17").Value = 0.4298245614035088
$ws.Range("K17").Value = 0.09210526315789473
$ws.Range("M17").Value = 0.01754385964912281
$ws.Range("O17").Value = 0.06140350877192982
$ws.Range("S17").Value = 0.08991228070175439
$ws.Range("F18").Value = 0.009216589861751152
$ws.Range("H18").Value = 0.2073732718894009
$ws.Range("I18").Value = 0.09216589861751152
$ws.Range("J18").Value = 0.3778801843317972
$ws.Range("K18").Value = 0.1152073732718894
$ws.Range("M18").Value = 0.02304147465437788
$ws.Range("N18").Value = 0.004608294930875576
$ws.Range("O18").Value = 0.05069124423963134
$ws.Range("S18").Value = 0.119815668202765
$ws.Range("F19").Value = 0.01487862176977291
$ws.Range("H19").Value = 0.2286609240407204
$ws.Range("I19").Value = 0.07674236491777604
$ws.Range("J19").Value = 0.350822239624119
$ws.Range("K19").Value = 0.1323414252153485
$ws.Range("M19").Value = 0.01018010963194988
$ws.Range("N19").Value = 0.0007830853563038371
$ws.Range("O19").Value = 0.07282693813625685
$ws.Range("S19").Value = 0.1127642913077525
